$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-30 Sunday", "2025-03-31 Monday"),
    @("646÷4=", "892÷2="),
    @("525÷8=", "825÷9="),
    @("195÷6=", "798÷7="),
    @("570÷7=", "994÷6="),
    @("747÷3=", "966÷8="),
    @("613÷2=", "284÷8="),
    @("394÷2=", "382÷2="),
    @("704÷9=", "823÷7="),
    @("233÷2=", "291÷2="),
    @("413÷2=", "548÷6="),
    @("238÷6=", "484÷5="),
    @("993÷8=", "101÷8="),
    @("217÷2=", "858÷2="),
    @("859÷6=", "374÷6="),
    @("436÷5=", "567÷9="),
    @("898÷4=", "100÷5="),
    @("627÷9=", "866÷9="),
    @("989÷8=", "520÷8="),
    @("117÷2=", "415÷8="),
    @("109÷4=", "752÷5="),
    @("105÷9=", "362÷9="),
    @("432÷4=", "449÷3="),
    @("919÷9=", "505÷8="),
    @("847÷8=", "442÷5="),
    @("758÷6=", "388÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
